# Applies the crypto price/symbol list update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.70"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.287"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05763"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.473"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.149"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8163"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8527"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1357"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06950"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03137"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02942"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09391"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.745"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04647"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005978"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006112"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004611"
$ws.Range("E21").Value = "20HotbitTokenHTBBestin24h"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00006105"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.498"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.137"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3196"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002335"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03669"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006255"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1053"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002860"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008506"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005260"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3705"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002301"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"

Write-Output "Applied 42 cell updates."
